# B1/B2 PowerPoint update:
#  1. Table on slide 5 gets a new table style id.
#  2. The deck's theme palette is switched from the "Integral" (Red
#     Violet) design to the plain "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the financial-documents table on slide 5 -----------------
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shape = $slide5.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{D93CB9CD-95A4-4322-B4BB-C708CA30062D}")
    }
}

# --- 2. Swap the theme colour palette to the Office defaults --------------
# Office theme colours (RRGGBB) in ThemeColorScheme.Item() order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # COM RGB values are packed as 0x00BBGGRR
    $bgr = ($b * 65536) + ($g * 256) + $r
    $themeColors.Item($i).RGB = $bgr
}
